$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.615425109863281
$ws.Range("B1").Value = 2.410564184188843
$ws.Range("C1").Value = 1.797758102416992
$ws.Range("D1").Value = 1.688449263572693
$ws.Range("E1").Value = 1.770780920982361
